$wb = $excel.ActiveWorkbook

# Trade #10 (base_strategy, DOWN) opened at 2026-02-16 22:58:16, mirrored into
# both the "All Trades" log and the per-strategy "base_strategy" sheet as a
# new row appended right after the existing last trade row (row 10).
#
# We copy the previous trade row first (instead of writing every cell value
# from scratch) so that cells holding date-/time-looking literal text (e.g.
# "2026-02-16") keep their original text representation instead of being
# auto-converted into an Excel date serial number, and so blank cells keep
# matching. Afterwards we only overwrite the handful of cells that actually
# differ for the new trade: Trade # and Time.

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $srcRow = 10
    $newRow = 11

    $ws.Range("A$srcRow`:Q$srcRow").Copy($ws.Range("A$newRow`:Q$newRow"))

    $ws.Cells.Item($newRow, 1).Value = 10
    $ws.Cells.Item($newRow, 3).Value = "22:58:16"
}
